# Fruta / hortaliza, semanal
# Insert a new week of price data (date 2021-11-08, serial 44508) for
# "Palta" / "Hass" at the top of the data block (row 388), pushing the
# existing rows down by three positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the three new records; everything at/after row 388 shifts
# down to rows 391+ (dimension grows from T493 to T496 automatically).
$ws.Rows("388:390").Insert()

# Row 388: Hass / Especial
$ws.Cells.Item(388,1).Value = 8
$ws.Cells.Item(388,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(388,3).Value = "Coquimbo"
$ws.Cells.Item(388,4).Value = 44508
$ws.Cells.Item(388,5).Value = 4
$ws.Cells.Item(388,6).Value = "Fruta"
$ws.Cells.Item(388,7).Value = 100106
$ws.Cells.Item(388,8).Value = "Oleaginosos"
$ws.Cells.Item(388,9).Value = 100106002
$ws.Cells.Item(388,10).Value = "Palta"
$ws.Cells.Item(388,11).Value = "Hass"
$ws.Cells.Item(388,12).Value = "Especial"
$ws.Cells.Item(388,13).Value = 200
$ws.Cells.Item(388,14).Value = 2400
$ws.Cells.Item(388,15).Value = 2500
$ws.Cells.Item(388,16).Value = 2450
$ws.Cells.Item(388,17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(388,18).Value = "Provincia de Limarí"
$ws.Cells.Item(388,19).Value = 2450
$ws.Cells.Item(388,20).Value = 1

# Row 389: Hass / Primera
$ws.Cells.Item(389,1).Value = 8
$ws.Cells.Item(389,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(389,3).Value = "Coquimbo"
$ws.Cells.Item(389,4).Value = 44508
$ws.Cells.Item(389,5).Value = 4
$ws.Cells.Item(389,6).Value = "Fruta"
$ws.Cells.Item(389,7).Value = 100106
$ws.Cells.Item(389,8).Value = "Oleaginosos"
$ws.Cells.Item(389,9).Value = 100106002
$ws.Cells.Item(389,10).Value = "Palta"
$ws.Cells.Item(389,11).Value = "Hass"
$ws.Cells.Item(389,12).Value = "Primera"
$ws.Cells.Item(389,13).Value = 400
$ws.Cells.Item(389,14).Value = 2100
$ws.Cells.Item(389,15).Value = 2200
$ws.Cells.Item(389,16).Value = 2150
$ws.Cells.Item(389,17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(389,18).Value = "Provincia de Limarí"
$ws.Cells.Item(389,19).Value = 2150
$ws.Cells.Item(389,20).Value = 1

# Row 390: Hass / Segunda
$ws.Cells.Item(390,1).Value = 8
$ws.Cells.Item(390,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(390,3).Value = "Coquimbo"
$ws.Cells.Item(390,4).Value = 44508
$ws.Cells.Item(390,5).Value = 4
$ws.Cells.Item(390,6).Value = "Fruta"
$ws.Cells.Item(390,7).Value = 100106
$ws.Cells.Item(390,8).Value = "Oleaginosos"
$ws.Cells.Item(390,9).Value = 100106002
$ws.Cells.Item(390,10).Value = "Palta"
$ws.Cells.Item(390,11).Value = "Hass"
$ws.Cells.Item(390,12).Value = "Segunda"
$ws.Cells.Item(390,13).Value = 340
$ws.Cells.Item(390,14).Value = 1900
$ws.Cells.Item(390,15).Value = 2000
$ws.Cells.Item(390,16).Value = 1950
$ws.Cells.Item(390,17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(390,18).Value = "Provincia de Limarí"
$ws.Cells.Item(390,19).Value = 1950
$ws.Cells.Item(390,20).Value = 1
